$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number need the
# NumberFormat="@" (Text) trick so Excel stores them as text, not a number,
# matching the source inlineStr cells. Style is reset to Normal afterward
# so the cell keeps its original (default) style index.
$ws.Range("D2").Value = "37.171.96"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.002.00"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.624"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.09%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.56%  "
$ws.Range("D14").Value = "2.295.81"
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.841"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").Value = "1.996.42"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "37.120.84"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.141"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("E30").Value = "  +11.56%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0658"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.68%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("E38").Value = "  -7.78%  "
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0214"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.18%  "
$ws.Range("D46").Value = "1.368.43"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +11.79%  "
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.39%  "
